$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.946.24'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.633.81'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''216.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").Value = '''0.5116'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''0.2579'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("D9").Value = '''0.06352'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '''19.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").Value = '''0.07784'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '''4.284'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '1.636.59'
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("D14").Value = '1.858.65'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '0.5524'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '''63.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '0.0₅7655'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '25.963.95'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '''195.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = '''4.428'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '''9.873'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = '6.048'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '1.891'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").Value = '141.84'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").Value = '''0.1260'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("D28").Value = '''15.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").Value = '''6.760'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").Value = '1.242'
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").Value = '''0.04887'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = '''3.245'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '''3.196'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").Value = '''1.543'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("D35").Value = '2.371'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '''0.8985'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").Value = '''0.5520'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.79%  '
$ws.Range("D38").Value = '''2.540'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("D39").Value = '1.118.67'
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").Value = '''0.01559'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = '''1.000'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").Value = '''5.595'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.69%  '
$ws.Range("D43").Value = '0.7974'
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").Value = '''97.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").Value = '1.770.11'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").Value = '0.0₈116'
$ws.Range("E46").Value = '  -7.46%  '
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").Value = '''54.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("D51").Value = '''7.568'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.10%  '
